$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: name + score + paidFees flag
$ws.Range("A2").Value = "Loveth"
$ws.Range("B2").Value = 37
$ws.Range("D2").Value = $true

# Update row 3: name + score
$ws.Range("A3").Value = "Excel"
$ws.Range("B3").Value = 73

# Move active selection as recorded in the saved workbook
[void]$ws.Range("O5").Select()
